$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 130803173
$ws.Range("B5").Value = 57884
$ws.Range("Q5").Value = 466010
$ws.Range("R5").Value = 7072292
$ws.Range("AC5").Value = "Ringhack äldre"
$ws.Range("A6").Value = 130803168
$ws.Range("B6").Value = 57884
$ws.Range("Q6").Value = 465993
$ws.Range("R6").Value = 7072380
$ws.Range("AC6").Value = "Ringhack färska och äldre"
$ws.Range("A7").Value = 130803170
$ws.Range("B7").Value = 57884
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("Q7").Value = 465964
$ws.Range("R7").Value = 7072357
$ws.Range("AC7").Value = "Ringhack färska och äldre"
$ws.Range("A8").Value = 130803176
$ws.Range("B8").Value = 91828
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H8").ClearContents()
$ws.Range("Q8").Value = 466028
$ws.Range("R8").Value = 7072541
$ws.Range("AC8").ClearContents()
$ws.Range("B9").Value = 79243
$ws.Range("B10").Value = 79243
$ws.Range("B11").Value = 91804
$ws.Range("B12").Value = 91828
$ws.Range("A13").Value = 130803165
$ws.Range("B13").Value = 92530
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 3298
$ws.Range("F13").Value = "Trådticka"
$ws.Range("G13").Value = "Climacocystis borealis"
$ws.Range("H13").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q13").Value = 466019
$ws.Range("R13").Value = 7072412
$ws.Range("AC13").ClearContents()
$ws.Range("B14").Value = 57884
$ws.Range("A15").Value = 130803166
$ws.Range("B15").Value = 57884
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("Q15").Value = 466057
$ws.Range("R15").Value = 7072377
$ws.Range("AC15").Value = "Ringhack"
$ws.Range("B16").Value = 57884
$ws.Range("B17").Value = 91828
$ws.Range("A18").Value = 130803169
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("Q18").Value = 465965
$ws.Range("R18").Value = 7072357
$ws.Range("AC18").Value = "Ringhack"
$ws.Range("A19").Value = 130803177
$ws.Range("B19").Value = 91828
$ws.Range("E19").Value = 5432
$ws.Range("F19").Value = "Granticka"
$ws.Range("G19").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H19").ClearContents()
$ws.Range("Q19").Value = 465951
$ws.Range("R19").Value = 7072435
$ws.Range("AC19").ClearContents()
